# Apply BOM updates: removed R49, removed C117 and fixed C116 and J2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Board - Rev A")

# --- Row 14: C116 (was C116, C117) - fixed part, now C0805 / 250V 1000pF cap ---
$ws.Range("B14").Value = "C116"
$ws.Range("H14").Value = "445-2277-1-ND"
$ws.Range("F14").Value = "C2012X7R2E102K"
$ws.Range("I14").Value = "CAP CER 1000PF 250V 10% X7R 0805"
$ws.Range("D14").Value = "C0805"
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = 0.17

# --- Row 17: J2 - RJ45 connector replaced with Pulse Electronics pulsejack part ---
$ws.Range("H17").Value = "553-1485-ND"
$ws.Range("F17").Value = "J0011D21BNL"
$ws.Range("I17").Value = "CONN PULSEJACK 1PORT 10/100B-TX"
$ws.Range("C17").Value = "CON-RJ45-J0011D21BNL"
$ws.Range("D17").Value = "CON-RJ45-J0011D21BNL"
$ws.Range("K17").Value = 7.1

# --- Row 33: R0402 0-ohm group - removed R49 ---
$ws.Range("B33").Value = "R4, R11, R12, R13, R15, R16, R17, R18, R74, R75"
$ws.Range("J33").Value = 11

# Recalculate dependent formulas
$excel.Calculate() | Out-Null

# Restore the cell-selection view state seen after the edit
$ws.Range("J34").Select() | Out-Null
